$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Apply each cell update. D/E columns hold numeric-looking text (prices,
# percentages) that must stay literal strings, so force Text format first.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '295.09'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '-4.40%'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '40.17'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '-2.12%'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '5.039'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '-3.24%'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.07422'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '-3.51%'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '4.330'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '0.68%'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.581'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '-4.01%'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.9252'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '1.25%'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.1168'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '-5.37%'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.1744'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '-3.80%'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08761'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '-3.98%'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.04181'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '-1.44%'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.1055'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '0.35%'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.001266'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '0.79%'
$ws.Range('B16').Value = 'CoinExToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.03859'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '-4.01%'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.005973'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '2.56%'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.363'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '0.56%'
$ws.Range('B19').Value = 'BitpandaEcosystemToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.3346'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '0.33%'
$ws.Range('B20').Value = 'MCDex'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.617'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '4.34%'
$ws.Range('B21').Value = 'ProBitToken'
$ws.Range('C21').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.1357'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '-1.75%'
$ws.Range('B22').Value = 'ZBToken'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.2816'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '4.00%'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '2.27%'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.003635'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '-14.78%'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.0001308'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '0.70%'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0003742'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '-95.02%'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02326'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '-7.53%'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05007'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '-5.58%'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.007744'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '-1.29%'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '124.38%'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '11.54%'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.007247'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '-9.21%'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.3210'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '4.58%'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.00006475'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '-3.48%'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '0.69%'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '-27.22%'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.004223'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '36.34%'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.00002113'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '0.69%'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '0.69%'
